# Applies the cryptocurrency price/volume updates described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "41.949.62"
$ws.Cells.Item(2, 5).Value = "  -2.33%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.291.28"
$ws.Cells.Item(3, 5).Value = "  -3.13%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.02%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'317.06"
$ws.Cells.Item(5, 5).Value = "  -0.29%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'103.10"
$ws.Cells.Item(6, 5).Value = "  -4.96%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "'0.624"
$ws.Cells.Item(7, 5).Value = "  -2.18%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.01%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.604"
$ws.Cells.Item(9, 5).Value = "  -3.58%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'39.23"
$ws.Cells.Item(10, 5).Value = "  -5.29%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -2.73%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  -4.24%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  -0.85%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'0.958"
$ws.Cells.Item(14, 5).Value = "  -5.44%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "'15.20"
$ws.Cells.Item(15, 5).Value = "  -5.08%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "2.639.41"
$ws.Cells.Item(16, 5).Value = "  -3.11%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "2.299.86"
$ws.Cells.Item(17, 5).Value = "  -2.74%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "42.028.25"
$ws.Cells.Item(18, 5).Value = "  -2.12%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "'7.36"
$ws.Cells.Item(19, 5).Value = "  -3.66%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "'0.0000106"
$ws.Cells.Item(20, 5).Value = "  -1.06%  "

# Row 21
$ws.Cells.Item(21, 2).Value = "Litecoin"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(21, 4).Value = "'73.25"
$ws.Cells.Item(21, 5).Value = "  -4.20%  "

# Row 22
$ws.Cells.Item(22, 2).Value = "PancakeSwap"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(22, 4).Value = "'3.61"
$ws.Cells.Item(22, 5).Value = "  -0.86%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "'279.05"
$ws.Cells.Item(23, 5).Value = "  +4.27%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'10.18"
$ws.Cells.Item(24, 5).Value = "  +7.49%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  -2.97%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +0.70%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'10.76"
$ws.Cells.Item(27, 5).Value = "  -6.41%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  +4.75%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'22.74"
$ws.Cells.Item(29, 5).Value = "  -3.18%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "'35.85"
$ws.Cells.Item(30, 5).Value = "  -3.05%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'163.35"
$ws.Cells.Item(31, 5).Value = "  -2.91%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'0.0870"
$ws.Cells.Item(32, 5).Value = "  -4.51%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "'2.83"
$ws.Cells.Item(33, 5).Value = "  -2.72%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "'5.78"
$ws.Cells.Item(34, 5).Value = "  -3.87%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  +3.22%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -5.14%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  -5.84%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  -4.98%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "'2.82"
$ws.Cells.Item(39, 5).Value = "  +2.94%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "'3.74"
$ws.Cells.Item(40, 5).Value = "  -3.58%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "'98.92"
$ws.Cells.Item(41, 5).Value = "  -7.00%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  -5.20%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "'69.13"
$ws.Cells.Item(43, 5).Value = "  -3.21%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(44, 4).Value = "'1.00"
$ws.Cells.Item(44, 5).Value = "  +0.16%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "Algorand"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(45, 4).Value = "'0.224"
$ws.Cells.Item(45, 5).Value = "  -6.34%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "'11.89"
$ws.Cells.Item(46, 5).Value = "  -4.56%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "'112.37"
$ws.Cells.Item(47, 5).Value = "  -1.44%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "'76.72"
$ws.Cells.Item(48, 5).Value = "  -0.10%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  -3.06%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'5.27"
$ws.Cells.Item(50, 5).Value = "  -5.29%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "1.587.08"
$ws.Cells.Item(51, 5).Value = "  +0.61%  "
